# Edit script: applies the diff changes to database_doc__sd_.docx
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1. Doctors table paragraph - extend sentence about speciality/patients
Replace-Text `
    "The primary key of this table is doctorID. Names, sex, speciality and date of birth are stored so the system has some information about the doctors and the speciality is especially important as the staff will need to be able to assign the right doctors to the patients." `
    "The primary key of this table is doctorID. Names, sex, speciality and date of birth are stored so the system has some information about the doctors the same as the patients but the specialty is especially important as the patients will need to be able to choose the right doctors for them, and the staff will be able to choose the right doctors for patients if necessary."

# 2. Bookings table paragraph - replace "It also needs..." sentence with doctor sentence
Replace-Text `
    ". It also needs to be stored by the system so it can be verified that the patient has made a booking and so the doctor can check on their past and future bookings " `
    ". As for the doctor, they will need it to check on their upcoming bookings and to make sure they don’t get double booked. "

# 3. Messages table paragraph - append sentence about pID relationship
Replace-Text `
    "The message body is used to send patients a confirmation message when they book doctor visits, change their doctors, etc." `
    "The message body is used to send patients a confirmation message when they book doctor visits, change their doctors, etc. with the pID forming that relationship and ensuring the message is sent to the right patient."

# 4. Patient Access/History table paragraph - full rewrite
Replace-Text `
    "The patientID is the both the primary and foreign key as it is the only ID of the table but has holds information from the Patients table, giving both tables a relationship. The most important attributes are dateAccessed and functionality as they are needed to track when and what the patient has accessed from the system." `
    "The primary key is a composite of all data types in this table, the PatientID, dateAccessed, and Functionality. This is because this is the only way to guarantee a unique entry into this table, as any given functionality can be accessed multiple times on different dates. The patientID also acts as the foreign key as it is the only ID of the table but also because it holds information from the Patients table, giving both tables the needed relationship. The most important attributes are dateAccessed and functionality as they are needed to track when and what the patient has accessed from the system."

# 5. Visit details / prescriptions paragraph - several small edits
Replace-Text `
    "are both the primary keys and foreign keys of the table as the information in the table is relevant to the Patient, Doctor and Booking tables. The visitDetails attribute holds information most likely from the doctor about any relevant details from their patients visits for the" `
    "are both the primary keys (composite) and foreign keys of the table as the information in this table is relevant to the Patient, Doctor and Booking tables. The visitDetails attribute holds information from the doctor about any relevant details from their patients visits for the"

Replace-Text "[patients]" "(patients)"

Replace-Text `
    " prescriptions given by their doctor." `
    " prescriptions given by their doctor, which is important if this was mentioned in the visit details as the result of the meeting or just an ongoing treatment."
